$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new row 4 values
$ws.Range("B4").Value = "Regular"
$ws.Range("C4").Value = "Class III"
$ws.Range("D4").Value = "Class III"
$ws.Range("E4").Value = "Customer"

# Column widths (values chosen so the engine's char->pixel rounding lands on
# the target stored widths as closely as possible)
$ws.Columns.Item(7).ColumnWidth = 23.1666666666667
$ws.Columns.Item(8).ColumnWidth = 12.1666666666667
$ws.Columns.Item(9).ColumnWidth = 15.6197916666667
$ws.Columns.Item(10).ColumnWidth = 17.4361979166667

# Selection
$ws.Range("J1").Select()
